$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf1"
$ws.Range("C2").Value = "Fgfr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 2.004760666666666
$ws.Range("H2").Value = 6.014282
$ws.Range("I2").Value = 0.1200698528618338
$ws.Range("J2").Value = 0.1200698528618338
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.4871643333333333
$ws.Range("N2").Value = 1.461493
$ws.Range("O2").Value = 0.1914458468300136
$ws.Range("P2").Value = 0.1914458468300136
$ws.Range("Q2").Value = 0.9766478936695554
$ws.Range("R2").Value = 8.789831043025998
$ws.Range("S2").Value = 0.0229868746598889
$ws.Range("T2").Value = 0.0229868746598889

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf1"
$ws.Range("C3").Value = "Fgfr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 2.004760666666666
$ws.Range("H3").Value = 6.014282
$ws.Range("I3").Value = 0.1200698528618338
$ws.Range("J3").Value = 0.1200698528618338
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.864751
$ws.Range("N3").Value = 5.594253
$ws.Range("O3").Value = 0.7328098752209857
$ws.Range("P3").Value = 0.7328098752209857
$ws.Range("Q3").Value = 3.738379457927333
$ws.Range("R3").Value = 33.645415121346
$ws.Range("S3").Value = 0.08798837389348252
$ws.Range("T3").Value = 0.08798837389348252

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf1"
$ws.Range("C4").Value = "Fgfr2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 2.004760666666666
$ws.Range("H4").Value = 6.014282
$ws.Range("I4").Value = 0.1200698528618338
$ws.Range("J4").Value = 0.1200698528618338
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1927433333333333
$ws.Range("N4").Value = 0.57823
$ws.Range("O4").Value = 0.07574427794900063
$ws.Range("P4").Value = 0.07574427794900063
$ws.Range("Q4").Value = 0.3864042534288889
$ws.Range("R4").Value = 3.47763828086
$ws.Range("S4").Value = 0.009094604308462345
$ws.Range("T4").Value = 0.009094604308462345

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fgf1"
$ws.Range("C5").Value = "Fgfr2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 9.409654999999999
$ws.Range("H5").Value = 28.228965
$ws.Range("I5").Value = 0.5635664696121425
$ws.Range("J5").Value = 0.5635664696121424
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.4871643333333333
$ws.Range("N5").Value = 1.461493
$ws.Range("O5").Value = 0.1914458468300136
$ws.Range("P5").Value = 0.1914458468300136
$ws.Range("Q5").Value = 4.584048304971666
$ws.Range("R5").Value = 41.25643474474499
$ws.Range("S5").Value = 0.1078924600198977
$ws.Range("T5").Value = 0.1078924600198977

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fgf1"
$ws.Range("C6").Value = "Fgfr2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 9.409654999999999
$ws.Range("H6").Value = 28.228965
$ws.Range("I6").Value = 0.5635664696121425
$ws.Range("J6").Value = 0.5635664696121424
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.864751
$ws.Range("N6").Value = 5.594253
$ws.Range("O6").Value = 0.7328098752209857
$ws.Range("P6").Value = 0.7328098752209857
$ws.Range("Q6").Value = 17.546663570905
$ws.Range("R6").Value = 157.919972138145
$ws.Range("S6").Value = 0.4129870742752056
$ws.Range("T6").Value = 0.4129870742752055

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgf1"
$ws.Range("C7").Value = "Fgfr2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 9.409654999999999
$ws.Range("H7").Value = 28.228965
$ws.Range("I7").Value = 0.5635664696121425
$ws.Range("J7").Value = 0.5635664696121424
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.1927433333333333
$ws.Range("N7").Value = 0.57823
$ws.Range("O7").Value = 0.07574427794900063
$ws.Range("P7").Value = 0.07574427794900063
$ws.Range("Q7").Value = 1.813648270216667
$ws.Range("R7").Value = 16.32283443195
$ws.Range("S7").Value = 0.04268693531703914
$ws.Range("T7").Value = 0.04268693531703913

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Fgf1"
$ws.Range("C8").Value = "Fgfr2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.282203999999999
$ws.Range("H8").Value = 15.846612
$ws.Range("I8").Value = 0.3163636775260238
$ws.Range("J8").Value = 0.3163636775260238
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.4871643333333333
$ws.Range("N8").Value = 1.461493
$ws.Range("O8").Value = 0.1914458468300136
$ws.Range("P8").Value = 0.1914458468300136
$ws.Range("Q8").Value = 2.573301390190666
$ws.Range("R8").Value = 23.159712511716
$ws.Range("S8").Value = 0.06056651215022697
$ws.Range("T8").Value = 0.06056651215022697

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Fgf1"
$ws.Range("C9").Value = "Fgfr2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.282203999999999
$ws.Range("H9").Value = 15.846612
$ws.Range("I9").Value = 0.3163636775260238
$ws.Range("J9").Value = 0.3163636775260238
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.864751
$ws.Range("N9").Value = 5.594253
$ws.Range("O9").Value = 0.7328098752209857
$ws.Range("P9").Value = 0.7328098752209857
$ws.Range("Q9").Value = 9.849995191203998
$ws.Range("R9").Value = 88.64995672083599
$ws.Range("S9").Value = 0.2318344270522976
$ws.Range("T9").Value = 0.2318344270522976

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Fgf1"
$ws.Range("C10").Value = "Fgfr2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 5.282203999999999
$ws.Range("H10").Value = 15.846612
$ws.Range("I10").Value = 0.3163636775260238
$ws.Range("J10").Value = 0.3163636775260238
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.1927433333333333
$ws.Range("N10").Value = 0.57823
$ws.Range("O10").Value = 0.07574427794900063
$ws.Range("P10").Value = 0.07574427794900063
$ws.Range("Q10").Value = 1.018109606306667
$ws.Range("R10").Value = 9.162986456759999
$ws.Range("S10").Value = 0.02396273832349915
$ws.Range("T10").Value = 0.02396273832349915
